# Insert two new rows for input/output tensor name parameters right above
# the existing "mlp_structure" block (old row 27, now row 29), shifting the
# rest of the "Specifications (structure)" table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("config_file_structure")

$ws.Activate()

$ws.Rows("27:28").Insert()

# New parameter-block row (input tensor name)
$ws.Range("E27").Value = '"input_tensor_name": <parameter block>,'
# New parameter-block row (output tensor name)
$ws.Range("E28").Value = '"output_tensor_name": <parameter block>,'

# Matching comments in column T
$ws.Range("T27").Value = "// input tensor name"
$ws.Range("T28").Value = "// output tensor name"

# Match the author's final selection/view state
$ws.Range("T29").Select()
